$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 9
$ws1.Range("F5").Value = 54
$ws1.Range("F6").Value = 870
$ws1.Range("F8").Value = 6901
$ws1.Range("F11").Value = 146
$ws1.Range("F12").Value = 6464
$ws1.Range("F15").Value = 4391
$ws1.Range("F19").Value = 4392
$ws1.Range("F21").Value = 237
$ws1.Range("F22").Value = 241
$ws1.Range("F23").Value = 327
$ws1.Range("F25").Value = 126
$ws1.Range("F26").Value = 168
$ws1.Range("F27").Value = 41
$ws1.Range("F29").Value = 74
$ws1.Range("F30").Value = 7947
$ws1.Range("F31").Value = 51
$ws1.Range("F32").Value = 1355
$ws1.Range("F33").Value = 665
$ws1.Range("F38").Value = 1604
$ws1.Range("F40").Value = 925
$ws1.Range("F42").Value = 4021
$ws1.Range("F45").Value = 109
$ws1.Range("F46").Value = 42
$ws1.Range("F48").Value = 1093
$ws1.Range("F49").Value = 7

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 13
$ws2.Range("F18").Value = 3
$ws2.Range("F19").Value = 872

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 9
$ws4.Range("F8").Value = 54
$ws4.Range("F9").Value = 13
$ws4.Range("F10").Value = 870
$ws4.Range("F12").Value = 6901
$ws4.Range("F15").Value = 146
$ws4.Range("F16").Value = 6464
$ws4.Range("F19").Value = 4391
$ws4.Range("F22").Value = 4392
$ws4.Range("F24").Value = 237
$ws4.Range("F25").Value = 241
$ws4.Range("F26").Value = 327
$ws4.Range("F28").Value = 126
$ws4.Range("F30").Value = 7947
$ws4.Range("F31").Value = 51
$ws4.Range("F32").Value = 1355
$ws4.Range("F33").Value = 665
$ws4.Range("F38").Value = 1604
$ws4.Range("F40").Value = 925
$ws4.Range("F42").Value = 4021
$ws4.Range("F45").Value = 109
$ws4.Range("F47").Value = 1093
$ws4.Range("F49").Value = 7

$wb.Save()
